# Update "Training Dashboard" sheet: recompute PERIOD TO EXPIRE (col H)
# and LAST UPDATE (col I) for the progress date of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Scratch cell, well outside the used range (A1:K34), used only to stage the
# new "LAST UPDATE" text so it can be copied in as a plain value. Writing the
# date-looking string straight into a General-formatted target cell makes
# Excel auto-convert it to a real date serial; staging it on a cell that is
# pre-formatted as Text and then pasting VALUES ONLY into the destination
# keeps the destination cell's own (unchanged) style/number format while
# still landing a literal text value.
$scratch = $ws.Cells.Item(100, 26)
$scratch.NumberFormat = "@"

for ($row = 3; $row -le 34; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    $currentPeriod = $hCell.Value2
    $hCell.Value = $currentPeriod - 1

    $scratch.Value = "04-Nov-2025"
    $scratch.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues - value only, keeps destination style
}

$scratch.Clear()
$excel.CutCopyMode = $false
